$wb = $excel.ActiveWorkbook
$rubric = $wb.Worksheets.Item("Rubric")
$grade = $wb.Worksheets.Item("Grade")

$grade.Range("A5").Value = "UX diagrams"
$rubric.Range("A5").Value = "UX Wireframes"

$grade.Activate()
$grade.Range("E12").Select()

$rubric.Activate()
$rubric.Range("A14").Select()
